$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add another person's survey review in row 15 (same pattern as the other
# entries: data rows are separated by one blank row).
$ws.Range("A15").Value = "Cam"
$ws.Range("B15").Value = "Tracking Macronutrients"
$ws.Range("C15").Value = "Carb Manager"
$ws.Range("D15").Value = "Tracks carbs, fat, protein. Able to track weight loss over weeks / months. Has a lot of brand names' macros"
$ws.Range("E15").Value = 'Always asks to get premium when opening app, cant get rid of "Carb" tracker bar'
$ws.Range("F15").Value = 'Ability to change "Carb" tracker bar into "Protein" tracker bar'
$ws.Range("G15").Value = 22
$ws.Range("H15").Value = "M"
$ws.Range("I15").Value = "Electrician"

# Columns C, D and E now hold longer text than before, so widen them to fit
# the new content (as Excel does automatically for bestFit columns).
$ws.Columns.Item(3).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(4).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(5).EntireColumn.AutoFit() | Out-Null

# Scroll the view back to the top-left and move the selection to just below
# the newly added row, matching where the author left off editing.
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("D16").Select() | Out-Null
